$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.885.00'
$ws.Range("E2").Value = '  -1.62%  '

$ws.Range("D3").Value = '2.456.66'
$ws.Range("E3").Value = '  -1.40%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '516.89'
$ws.Range("E5").Value = '  -2.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.27'
$ws.Range("E6").Value = '  -2.46%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  -1.77%  '

$ws.Range("D9").Value = '2.462.90'
$ws.Range("E9").Value = '  -1.66%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0971'
$ws.Range("E10").Value = '  -3.75%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.156'
$ws.Range("E11").Value = '  -0.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.28'
$ws.Range("E12").Value = '  -1.84%  '

$ws.Range("E13").Value = '  -3.37%  '

$ws.Range("D14").Value = '2.893.80'
$ws.Range("E14").Value = '  -1.49%  '

$ws.Range("D15").Value = '57.810.59'
$ws.Range("E15").Value = '  -1.58%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.86'
$ws.Range("E16").Value = '  -3.32%  '

$ws.Range("E17").Value = '  -2.78%  '

$ws.Range("D18").Value = '2.459.07'
$ws.Range("E18").Value = '  -1.80%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.54'
$ws.Range("E19").Value = '  -4.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '318.50'
$ws.Range("E20").Value = '  -1.00%  '

$ws.Range("E21").Value = '  -2.94%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.09%  '

$ws.Range("E23").Value = '  -4.95%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.41'
$ws.Range("E24").Value = '  -1.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.405'
$ws.Range("E25").Value = '  -3.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("E27").Value = '  -2.99%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.26'
$ws.Range("E28").Value = '  -2.85%  '

$ws.Range("D29").Value = '0.0₃0737'
$ws.Range("E29").Value = '  -2.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.02'
$ws.Range("E30").Value = '  -1.97%  '

$ws.Range("E31").Value = '  -3.39%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.19'
$ws.Range("E32").Value = '  -3.58%  '

$ws.Range("E33").Value = '  -0.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.04%  '

$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("E36").Value = '  -1.98%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.32'
$ws.Range("E37").Value = '  -1.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.92'
$ws.Range("E38").Value = '  -1.79%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.32'
$ws.Range("E39").Value = '  -1.06%  '

$ws.Range("E40").Value = '  -4.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.783'
$ws.Range("E41").Value = '  -0.59%  '

$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '270.17'
$ws.Range("E42").Value = '  -3.49%  '

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.40'
$ws.Range("E43").Value = '  -4.34%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.93'
$ws.Range("E44").Value = '  -1.46%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.586'
$ws.Range("E45").Value = '  -2.72%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '123.11'
$ws.Range("E46").Value = '  -4.99%  '

$ws.Range("E47").Value = '  -1.93%  '

$ws.Range("E48").Value = '  -2.96%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0210'
$ws.Range("E49").Value = '  -3.23%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.65'
$ws.Range("E50").Value = '  -2.64%  '

$ws.Range("D51").Value = '1.722.33'
$ws.Range("E51").Value = '  -1.54%  '
